# "add update tenant and profile url"
# Append two new test-plan rows to the Tenant Test Plan sheet:
#   Row 22: Tenant_Service_update_tenant_returns_successfully -> HTTP 200
#   Row 23: Tenant_Service_rejects_invalid_state               -> HTTP 400

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "Tenant_Service_update_tenant_returns_successfully"
$ws.Range("B22").Value = "HTTP 200"

$ws.Range("A23").Value = "Tenant_Service_rejects_invalid_state"
$ws.Range("B23").Value = "HTTP 400"

# Move the active selection to A23, matching the post-edit workbook state.
$ws.Range("A23").Select()
